$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark D-column (Price) cells as Text before assigning, to preserve exact
# formatted strings (e.g. trailing zeros / multi-dot separators) instead of
# Excel auto-converting them to numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D12", "D13", "D14", "D15", "D17", "D19", "D20", "D23", "D28", "D30", "D32", "D34", "D35", "D36", "D40", "D41", "D42", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price (D) updates
$ws.Range("D2").Value = "63.016.21"
$ws.Range("D3").Value = "2.552.24"
$ws.Range("D5").Value = "578.96"
$ws.Range("D6").Value = "147.15"
$ws.Range("D8").Value = "0.583"
$ws.Range("D12").Value = "0.354"
$ws.Range("D13").Value = "27.20"
$ws.Range("D14").Value = "3.009.65"
$ws.Range("D15").Value = "62.934.06"
$ws.Range("D17").Value = "2.549.67"
$ws.Range("D19").Value = "335.53"
$ws.Range("D20").Value = "4.34"
$ws.Range("D23").Value = "65.49"
$ws.Range("D28").Value = "8.35"
$ws.Range("D30").Value = "1.90"
$ws.Range("D32").Value = "177.64"
$ws.Range("D34").Value = "410.59"
$ws.Range("D35").Value = "19.12"
$ws.Range("D36").Value = "0.400"
$ws.Range("D40").Value = "1.00"
$ws.Range("D41").Value = "39.53"
$ws.Range("D42").Value = "151.42"
$ws.Range("D44").Value = "20.82"
$ws.Range("D46").Value = "0.603"
$ws.Range("D47").Value = "0.0967"
$ws.Range("D48").Value = "0.0238"
$ws.Range("D49").Value = "18.28"
$ws.Range("D51").Value = "11.31"

# Volume(1h) (E) updates — values already contain surrounding spaces so
# Excel keeps them as text automatically.
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("E29").Value = "  +3.28%  "
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("E31").Value = "  -3.00%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("E50").Value = "  -8.60%  "
$ws.Range("E51").Value = "  -0.04%  "
